$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replicate formatting for the new rows by copying existing formatted cells ---
# Standard data rows (time value in col A, Tahoma text in col B) -> rows 289-311 (23 rows)
$ws.Range("A265:B287").Copy() | Out-Null
$ws.Range("A289:B311").PasteSpecial(-4122) | Out-Null

# "Clase NN" section-header row (bold white text on dark fill) -> row 312
$ws.Range("A288:B288").Copy() | Out-Null
$ws.Range("A312:B312").PasteSpecial(-4122) | Out-Null

# --- Fill in the new cell values ---
$ws.Range("A289").Value = 0.0024305555555555556
$ws.Range("B289").Value = "Clientes de BD"
$ws.Range("A290").Value = 0.003125
$ws.Range("B290").Value = "Cliente GUI (Mongo Compass)"
$ws.Range("A291").Value = 0.003472222222222222
$ws.Range("B291").Value = "Cliente App (desde node ==> Mongoose)  // Cliente Web ==>  Mongo atlas"
$ws.Range("A292").Value = 0.004861111111111111
$ws.Range("B292").Value = "ODM ==> mongoose, es un traductor entre el codigo y la base de datos "
$ws.Range("A293").Value = 0.006944444444444444
$ws.Range("B293").Value = "DBaaS (Database as a Service) Mongo Altlas"
$ws.Range("A294").Value = 0.00798611111111111
$ws.Range("B294").Value = "Explicacion de modelos en la nube como servicio (Cloud Models)"
$ws.Range("A295").Value = 0.015625
$ws.Range("B295").Value = "Mongo Atlas"
$ws.Range("A296").Value = 0.016666666666666666
$ws.Range("B296").Value = "Configuracion de Mongo Atlas"
$ws.Range("A297").Value = 0.022222222222222223
$ws.Range("B297").Value = "Panel de Mongo Atlas"
$ws.Range("A298").Value = 0.022569444444444444
$ws.Range("B298").Value = "Configuracion Network Access IP"
$ws.Range("A299").Value = 0.023958333333333335
$ws.Range("B299").Value = "Configuracion de usuario"
$ws.Range("A300").Value = 0.025694444444444443
$ws.Range("B300").Value = "Conexion a base de datos"
$ws.Range("A301").Value = 0.029513888888888888
$ws.Range("B301").Value = "Conexión con MongoDB Compass"
$ws.Range("A302").Value = 0.034027777777777775
$ws.Range("B302").Value = "Mongoose ==> sirve para definir esquemas definidos"
$ws.Range("C302").Value = " "
$ws.Range("A303").Value = 0.03819444444444445
$ws.Range("B303").Value = "Conexión de cero, instalacion mongoose"
$ws.Range("A304").Value = 0.03993055555555555
$ws.Range("B304").Value = "donde poner el nombre de la `"base de datos`" en la url de conexión  "
$ws.Range("A305").Value = 0.04097222222222222
$ws.Range("B305").Value = "function de conexión a la base de datos"
$ws.Range("A306").Value = 0.043055555555555555
$ws.Range("B306").Value = "Creacion del Schema"
$ws.Range("A307").Value = 0.04652777777777778
$ws.Range("B307").Value = "Archivo registro.js"
$ws.Range("A308").Value = 0.05381944444444445
$ws.Range("B308").Value = "Metodos en CRUD en el product manager"
$ws.Range("A309").Value = 0.05694444444444444
$ws.Range("B309").Value = "Controllers"
$ws.Range("A310").Value = 0.059722222222222225
$ws.Range("B310").Value = "del manager se pasa al controller"
$ws.Range("A311").Value = 0.06805555555555555
$ws.Range("B311").Value = "repaso de las carpetas y archivos generados"

# Row 312: new "Clase 15" section header + first topic
$ws.Range("A312").Value = "Clase 15"
$ws.Range("B312").Value = "Primera práctica integradora"

# --- Column B width (diff: 106.140625 bestFit -> 115.7109375, bestFit removed) ---
$ws.Columns.Item(2).ColumnWidth = 115.7109375

# --- Sheet view: scroll position + new selection after the appended rows ---
$ws.Range("A313").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 280
